# Generate Report for Handoff
# The source file "b3da3f9c-c1fd-41ef-8dc5-0c5ee32270ec.md" (row 3 on each sheet) has been
# handed off for translation: update its status and stamp the handoff datetime for both
# the zh-cn and de-de locales, and reflect the new status on the Overview sheet.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: zh-cn / de-de status columns for the b3da3f9c row
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: status + latest handoff datetime for the b3da3f9c row
$wsZhCn.Range("B3").Value = "Ready for handoff"
$wsZhCn.Range("D3").Value = "2016-03-10 12:28:30"

# de-de sheet: status + latest handoff datetime for the b3da3f9c row
$wsDeDe.Range("B3").Value = "Ready for handoff"
$wsDeDe.Range("D3").Value = "2016-03-10 12:28:33"
